$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the dash in "As-designed" with a space -> "As designed"
$ws.Range("G4").Value = "As designed"
$ws.Range("G5").Value = "As designed"
$ws.Range("G7").Value = "As designed"

# Update the selected cell/range on the sheet view
$ws.Range("D11").Select()
